$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "240.41", "82.00") are preserved verbatim as text instead of
# being coerced into numbers, matching the original inline-string data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.301.62"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.841.38"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").Value = "0.9976"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").Value = "240.41"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "0.6284"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").Value = "0.9986"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "0.07451"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("D9").Value = "0.2893"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "24.34"
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("D11").Value = "0.07716"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.841.89"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "0.6781"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "0.00001014"
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("D16").Value = "82.00"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "6.138"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "29.363.93"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "228.35"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "0.9989"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "7.397"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "0.9996"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "158.64"
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "8.415"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("D27").Value = "17.54"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "0.06411"
$ws.Range("E28").Value = "  +14.27%  "
$ws.Range("D29").Value = "1.393"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "1.472"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "4.079"
$ws.Range("E31").Value = "  -1.36%  "
$ws.Range("D32").Value = "4.058"
$ws.Range("D33").Value = "1.821"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D35").Value = "0.6944"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").Value = "2.578"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "2.829"
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("D38").Value = "1.252.91"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").Value = "0.01810"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "6.521"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").Value = "0.9085"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "0.9978"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").Value = "2.005.19"
$ws.Range("E43").Value = "  -12.61%  "
$ws.Range("D44").Value = "101.15"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("D45").Value = "66.29"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").Value = "7.024"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "0.1168"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").Value = "8.989"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "1.673"
$ws.Range("E51").Value = "  -0.63%  "
